$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("M18").Value = -716
$ws.Range("H53").Value = 433.26666
$ws.Range("I53").Value = 532.7
$ws.Range("J53").Value = 234.4
$ws.Range("K53").Value = 532.7
$ws.Range("L53").Value = 234.4
$ws.Range("M53").Value = 104.3
$ws.Range("N53").Value = -1508.4
$ws.Range("H88").Value = 1281.5
$ws.Range("I88").Value = 425
$ws.Range("J88").Value = 1567
$ws.Range("K88").Value = 425
$ws.Range("L88").Value = 1567
$ws.Range("M88").Value = -19
$ws.Range("N88").Value = -2379
$ws.Range("H91").Value = 1281.5
$ws.Range("I91").Value = 425
$ws.Range("J91").Value = 1567
$ws.Range("K91").Value = 425
$ws.Range("L91").Value = 1567
$ws.Range("M91").Value = 979
$ws.Range("N91").Value = -4375
$ws.Range("H92").Value = 256.53845
$ws.Range("I92").Value = 271.22223
$ws.Range("J92").Value = 223.5
$ws.Range("K92").Value = 271.22223
$ws.Range("L92").Value = 223.5
$ws.Range("M92").Value = 976.7777699999999
$ws.Range("N92").Value = -2719.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H137").Value = 1196.9048
$ws.Range("I137").Value = 1081.5264
$ws.Range("J137").Value = 2293
$ws.Range("K137").Value = 3244.5792
$ws.Range("L137").Value = 6879
$ws.Range("M137").Value = -694.5792000000001
$ws.Range("N137").Value = -11979

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("K57").Value = 5000
$ws.Range("M57").Value = -4516
$ws.Range("H61").Value = 1873
$ws.Range("I61").Value = 1497.3334
$ws.Range("K61").Value = 1497.3334
$ws.Range("M61").Value = -1285.3334
$ws.Range("H63").Value = 1869.7273
$ws.Range("I63").Value = 1233.8
$ws.Range("J63").Value = 2399.6667
$ws.Range("K63").Value = 1233.8
$ws.Range("L63").Value = 2399.6667
$ws.Range("M63").Value = -547.8
$ws.Range("N63").Value = -3771.6667
$ws.Range("H66").Value = 1869.7273
$ws.Range("I66").Value = 1233.8
$ws.Range("J66").Value = 2399.6667
$ws.Range("K66").Value = 6169
$ws.Range("L66").Value = 11998.3335
$ws.Range("M66").Value = -2737
$ws.Range("N66").Value = -18862.3335
$ws.Range("H74").Value = 1490.1428
$ws.Range("I74").Value = 1490.1428
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1490.1428
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -616.1428000000001
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1490.1428
$ws.Range("I77").Value = 1490.1428
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7450.714
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -3082.714
$ws.Range("N77").ClearContents()
$ws.Range("H94").Value = 73582.25
$ws.Range("J94").Value = 73582.25
$ws.Range("L94").Value = 73582.25
$ws.Range("N94").Value = -75384.25
$ws.Range("H122").Value = 1837.4445
$ws.Range("I122").Value = 1837.4445
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5512.333500000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3062.333500000001
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2405.6428
$ws.Range("I132").Value = 2513.3845
$ws.Range("J132").Value = 1005
$ws.Range("K132").Value = 7540.1535
$ws.Range("L132").Value = 3015
$ws.Range("M132").Value = -5010.1535
$ws.Range("N132").Value = -8075
$ws.Range("H136").Value = 1873
$ws.Range("I136").Value = 1497.3334
$ws.Range("K136").Value = 4492.0002
$ws.Range("M136").Value = -1942.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2025.3182
$ws.Range("I86").Value = 1290.3125
$ws.Range("J86").Value = 3985.3333
$ws.Range("K86").Value = 1290.3125
$ws.Range("L86").Value = 3985.3333
$ws.Range("M86").Value = -167.3125
$ws.Range("N86").Value = -6231.3333
$ws.Range("H89").Value = 2025.3182
$ws.Range("I89").Value = 1290.3125
$ws.Range("J89").Value = 3985.3333
$ws.Range("K89").Value = 6451.5625
$ws.Range("L89").Value = 19926.6665
$ws.Range("M89").Value = -835.5625
$ws.Range("N89").Value = -31158.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2162.0417
$ws.Range("I16").Value = 1190.8667
$ws.Range("J16").Value = 3780.6667
$ws.Range("K16").Value = 1190.8667
$ws.Range("L16").Value = 3780.6667
$ws.Range("M16").Value = -903.8667
$ws.Range("N16").Value = -4354.6667
$ws.Range("H31").Value = 3981.889
$ws.Range("I31").Value = 4202.3335
$ws.Range("J31").Value = 3871.6667
$ws.Range("K31").Value = 4202.3335
$ws.Range("L31").Value = 3871.6667
$ws.Range("M31").Value = -3907.3335
$ws.Range("N31").Value = -4461.6667
$ws.Range("H34").Value = 3981.889
$ws.Range("I34").Value = 4202.3335
$ws.Range("J34").Value = 3871.6667
$ws.Range("K34").Value = 4202.3335
$ws.Range("L34").Value = 3871.6667
$ws.Range("M34").Value = -4000.3335
$ws.Range("N34").Value = -4275.6667
$ws.Range("H59").Value = 33936.75
$ws.Range("J59").Value = 33936.75
$ws.Range("L59").Value = 33936.75
$ws.Range("N59").Value = -36226.75
$ws.Range("H60").Value = 21333.334
$ws.Range("H69").Value = 4980.3335
$ws.Range("I69").Value = 4980.3335
$ws.Range("K69").Value = 4980.3335
$ws.Range("M69").Value = -4231.3335
$ws.Range("H70").Value = 17055.555
$ws.Range("J70").Value = 17055.555
$ws.Range("L70").Value = 17055.555
$ws.Range("N70").Value = -17685.555
$ws.Range("H72").Value = 4980.3335
$ws.Range("I72").Value = 4980.3335
$ws.Range("K72").Value = 14941.0005
$ws.Range("M72").Value = -11197.0005
$ws.Range("H73").Value = 17055.555
$ws.Range("J73").Value = 17055.555
$ws.Range("L73").Value = 17055.555
$ws.Range("N73").Value = -19239.555
$ws.Range("H86").Value = 5998.3335
$ws.Range("I86").Value = 6000
$ws.Range("K86").Value = 6000
$ws.Range("M86").Value = -4877
$ws.Range("H89").Value = 5998.3335
$ws.Range("I89").Value = 6000
$ws.Range("K89").Value = 30000
$ws.Range("M89").Value = -24384
$ws.Range("H99").Value = 4913.273
$ws.Range("I99").Value = 5583
$ws.Range("K99").Value = 5583
$ws.Range("M99").Value = -4085
$ws.Range("H105").Value = 1192.4286
$ws.Range("J105").Value = 900
$ws.Range("L105").Value = 900
$ws.Range("N105").Value = -4394
$ws.Range("H113").Value = 2162.0417
$ws.Range("I113").Value = 1190.8667
$ws.Range("J113").Value = 3780.6667
$ws.Range("K113").Value = 1190.8667
$ws.Range("L113").Value = 3780.6667
$ws.Range("M113").Value = 979.1333
$ws.Range("N113").Value = -8120.6667
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H126").Value = 4913.273
$ws.Range("I126").Value = 5583
$ws.Range("K126").Value = 16749
$ws.Range("M126").Value = -14279
$ws.Range("H134").Value = 7899.6
$ws.Range("I134").Value = 7110.6665
$ws.Range("K134").Value = 21331.9995
$ws.Range("M134").Value = -18796.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 765.13513
$ws.Range("I132").Value = 740.6774
$ws.Range("K132").Value = 6666.096600000001
$ws.Range("M132").Value = -4136.096600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 29999.8
$ws.Range("J45").Value = 29999.8
$ws.Range("L45").Value = 29999.8
$ws.Range("N45").Value = -31117.8
$ws.Range("H92").Value = 7754.8
$ws.Range("J92").Value = 7754.8
$ws.Range("L92").Value = 7754.8
$ws.Range("N92").Value = -11498.8
$ws.Range("H113").Value = 1899.75
$ws.Range("I113").Value = 1799.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1799.5
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 370.5
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 9619004
$ws.Range("I122").Value = 10420087
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 31260261
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -31257811
$ws.Range("N122").Value = -22900
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 2496
$ws.Range("I126").Value = 2496
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7488
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5018
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4351.8667
$ws.Range("I61").Value = 4529
$ws.Range("J61").Value = 3200.5
$ws.Range("K61").Value = 4529
$ws.Range("L61").Value = 3200.5
$ws.Range("M61").Value = -4327
$ws.Range("N61").Value = -3604.5
$ws.Range("H113").Value = 4351.8667
$ws.Range("I113").Value = 4529
$ws.Range("J113").Value = 3200.5
$ws.Range("K113").Value = 4529
$ws.Range("L113").Value = 3200.5
$ws.Range("M113").Value = -2359
$ws.Range("N113").Value = -7540.5
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H133").Value = 49999
$ws.Range("I133").Value = 49999
$ws.Range("K133").Value = 49999
$ws.Range("M133").Value = -47469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H41").Value = 15192.444
$ws.Range("J41").Value = 14285.25
$ws.Range("L41").Value = 14285.25
$ws.Range("N41").Value = -15065.25
$ws.Range("H62").Value = 2512.2856
$ws.Range("I62").Value = 2557.6
$ws.Range("K62").Value = 2557.6
$ws.Range("M62").Value = -1933.6
$ws.Range("H65").Value = 2512.2856
$ws.Range("I65").Value = 2557.6
$ws.Range("K65").Value = 12788
$ws.Range("M65").Value = -9668
$ws.Range("H136").Value = 1628.5625
$ws.Range("I136").Value = 1628.5625
$ws.Range("K136").Value = 4885.6875
$ws.Range("M136").Value = -2335.6875
